$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1) - order matters for shared-string indices
$ws.Range("C1").Value = "Is Business Owner/Company Employee"
$ws.Range("D1").Value = "Business/Company Name"
$ws.Range("E1").Value = "Business/Company Position"

# New data cells (row 2) - set E2 before D2 so shared strings land in the
# same order as the target workbook (Software Engineer before Template Co. Inc)
$ws.Range("C2").Value = $true
$ws.Range("E2").Value = "Software Engineer"
$ws.Range("D2").Value = "Template Co. Inc"

# Column widths for the new columns (closest achievable values through the
# pixel-snapping ColumnWidth setter)
$ws.Columns.Item(3).ColumnWidth = 33.666666666666664
$ws.Columns.Item(4).ColumnWidth = 39.333333333333336
$ws.Columns.Item(5).ColumnWidth = 36.5

# Match the selection left behind in the saved file
[void]$ws.Range("E3").Select()
